$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New skill row: ID 13 "LimiterRemoval" / StatusSkill, Mana 12, CoolDown 6
# Copy formatting from the row above (row 14) so the new row matches the
# existing table's look (fontId/fillId = "Bom" style), then fill values.
$ws.Range("A14:E14").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "LimiterRemoval"
$ws.Range("C15").Value = "StatusSkill"
$ws.Range("D15").Value = 12
$ws.Range("E15").Value = 6

# New empty formatted cell below the table (F16), styled like the existing
# "Ruim" (red) marker cell at H12 but using the red palette: underlined
# dark-red font on the red fill already inherited from the column style.
$ws.Range("F16").Font.Underline = 2
$ws.Range("F16").Font.Color = 393372

$ws.Range("F16").Select()
